$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.229832530021667
$ws.Range("B1").Value = 3.483012199401855
$ws.Range("C1").Value = 4.877372741699219
$ws.Range("D1").Value = 2.173002481460571
$ws.Range("E1").Value = 1.374904751777649
